# Add in better CDISC CT handling for specific fields
#
# The "study" sheet stored studyType/studyPhase as combined
# "CODE=Label" strings (e.g. "C98388=Interventional Study").
# Split these into separate label/code values so the two concepts
# are stored independently:
#   C2 (studyType)  "C98388=Interventional Study" -> "Interventional Study"
#   D2 (studyPhase) "C15602=Phase III Trial"       -> "C15602"

$wb = $excel.ActiveWorkbook

$wsStudy = $wb.Worksheets.Item("study")
$wsStudy.Range("C2").Value = "Interventional Study"
$wsStudy.Range("D2").Value = "C15602"

# The "study" tab becomes the active/selected tab (it was previously
# "studyIdentifiers"), with a new active cell selection. "studyIdentifiers"
# keeps its own last selection (F1) unchanged, it just stops being the
# active tab.
$wsStudy.Activate()
$wsStudy.Range("C14").Select()
